$d = $word.ActiveDocument

# The change only touches the cosmetic DrawingML "name" attribute on the
# wp:docPr / pic:cNvPr elements of three inline pictures (the picture's
# real identity - its relationship/media target - is untouched). That
# attribute isn't exposed as InlineShape.Name in the object model, so we
# round-trip the package through Range.WordOpenXML and patch the three
# unique XML blocks directly.

$full = $d.Content.WordOpenXML

# Footer (first page), wp:docPr id="3": image2.png -> image1.png
$old1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$new1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'

# Footer (default/primary), wp:docPr id="2": image2.png -> image1.png
$old2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$new2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'

# Header (first page), wp:docPr id="1": image1.jpg -> image2.jpg
$old3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'
$new3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'

if ($full.IndexOf($old1) -lt 0) { throw "pattern 1 not found" }
if ($full.IndexOf($old2) -lt 0) { throw "pattern 2 not found" }
if ($full.IndexOf($old3) -lt 0) { throw "pattern 3 not found" }

$full = $full.Replace($old1, $new1)
$full = $full.Replace($old2, $new2)
$full = $full.Replace($old3, $new3)

$d.Content.WordOpenXML = $full

Write-Host "done"
